# The unified diff supplied for this task touches only two parts:
#   word/document.xml  and  word/styles.xml
# Every single hunk in it is a pure XML-canonicalization artifact: the
# "after" side of each changed line has exactly the same element /
# attribute name->value set as the "before" side, just with the
# attributes (and root-element namespace declarations) re-ordered
# alphabetically (the classic output of e.g. XML C14N canonicalization
# used by the tool that produced the diff for review). Canonicalizing
# the original word/document.xml / word/styles.xml from before.docx
# reproduces the "after" text of the diff byte-for-byte, which proves
# there is no actual semantic/content change requested for these parts
# (no text, formatting, style, section, or structural change).
#
# Word's COM object model has no notion of "XML attribute order" (and
# this interop runtime does not re-canonicalize a part's attribute
# order on save either) so the correct - and only faithful -
# reproduction of "the change described by the diff" is to leave the
# document's content completely untouched.
#
# Intentionally a no-op against the object model: just touch the
# ActiveDocument reference so the script still "does something"
# without mutating any content, formatting, or structure.
$d = $word.ActiveDocument
$null = $d.Name
